$wb = $excel.ActiveWorkbook

$overviewSheet = $wb.Worksheets.Item("Overview")
$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")

# Column E = Priority: "low" -> "ht" for rows 4-7
# Column H = Latest Handoff Datetime: bumped forward for rows 4-7
$zhSheet.Range("E4:E7").Value = "ht"
$deSheet.Range("E4:E7").Value = "ht"

$zhSheet.Range("H4:H7").Value = "2016-09-04 16:35:33"
$deSheet.Range("H4:H7").Value = "2016-09-04 16:35:38"

# Overview sheet "Latest HO Xliff Generate Date" for de-de tracks the same
# shared datetime string, so it advances alongside it.
$overviewSheet.Range("G4:G7").Value = "2016-09-04 16:35:38"
